# Auto-generated: update cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.376.72'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.653.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.93'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.260'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0614'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.889.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.653.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.570'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.401.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -7.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0725'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.52%  '
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  -4.33%  '
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.459.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.907'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.571'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.84%  '
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.23%  '
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.795.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.71%  '
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.05%  '
